$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShape($shapeRange, $newName) {
    $count = $shapeRange.Count
    for ($i = 1; $i -le $count; $i++) {
        $ish = $shapeRange.Item($i)
        $shp = $ish.ConvertToShape()
        $shp.Name = $newName
        [void]$shp.ConvertToInlineShape()
    }
}

# Footer "default" (footer1.xml) - PearsonLogo image: image2.png -> image1.png
$ftrDefault = $sec.Footers.Item(1)
if ($ftrDefault.Exists) {
    Rename-InlineShape $ftrDefault.Range.InlineShapes "image1.png"
}

# Footer "first page" (footer2.xml) - PearsonLogo image: image2.png -> image1.png
$ftrFirst = $sec.Footers.Item(2)
if ($ftrFirst.Exists) {
    Rename-InlineShape $ftrFirst.Range.InlineShapes "image1.png"
}

# Header "first page" (header2.xml) - BTec_Logo-Orange image: image1.jpg -> image2.jpg
$hdrFirst = $sec.Headers.Item(2)
if ($hdrFirst.Exists) {
    Rename-InlineShape $hdrFirst.Range.InlineShapes "image2.jpg"
}
